$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new columns (G:K) for the "meta" statistics, shifting the existing
# arrecadado_sucesso..maior_ano columns from G:V to L:AA.
$ws.Columns("G:K").Insert(-4161)

# New header labels for the inserted columns.
$ws.Cells.Item(1,7).Value  = "meta"
$ws.Cells.Item(1,8).Value  = "meta_avg"
$ws.Cells.Item(1,9).Value  = "meta_std"
$ws.Cells.Item(1,10).Value = "meta_min"
$ws.Cells.Item(1,11).Value = "meta_max"

# Give the new columns the same "R$ #,##0.00" currency format used by the
# neighbouring arrecadado_* columns.
$ws.Range("G2:K22").NumberFormat = "R$ #,##0.00"

$rowData = @{
  2 = @(3429.901669220172,3429.901669220172,0,3429.901669220172,3429.901669220172,1906.113717404507,953.0568587022536,1131.80766455125,152.7479840991555,1753.365733305352,20.64400048281337,2.192741622159442,19.09349801239444,22.1945029532323,87,43.5,50.20458146424487,8,79,2018,2021)
  3 = @(4685.223515355117,1171.305878838779,1513.463401005405,203.3209614085382,3429.901669220172,2037.988779487729,407.5977558975457,752.9936319265861,40.66419228170764,1753.365733305352,19.42932575751352,6.3937578269931,12.5989877672314,27.95663219367401,95,19,33.57826678076163,2,79,2018,2023)
  4 = @(44668.21149582993,2233.410574791496,4669.826278860146,0,21176.91783511972,8971.61544582792,407.8007020830873,1065.395454170603,1.087396962410123,5087.076865717208,18.30903869114585,13.08741272121433,1.011042153300025,45.27801686072104,323,14.68181818181818,27.46146453909613,1,128,2016,2023)
  5 = @(5793.543450529515,1931.181150176505,2696.19764694338,203.3209614085382,5037.956787043533,6149.198283664386,1229.839656732877,2169.930045219226,6.627188424929332,5087.076865717208,24.17877060337857,14.1532578054747,6.627188424929332,39.74278801341569,166,33.2,53.54624169818084,1,128,2020,2023)
  6 = @(11142.9217614901,857.1478278069306,762.0287036244347,63.33892833870453,2420.445520432476,1963.34220364979,122.7088877281118,134.9786244833396,5.278244028225378,538.4389998789497,25.24042809340392,19.39110965726495,5.278244028225378,70.01644246718027,79,4.9375,4.234284670008226,1,15,2017,2023)
  7 = @(9396.888161248213,854.2625601134739,847.4430641318904,128.9715040670144,2829.876053795064,2039.453666840368,119.9678627553158,157.2645751483118,2.022084306600051,538.4389998789497,19.40166161997551,12.36470609657289,1.011042153300025,45.27801686072104,90,5.294117647058823,4.984534906618081,1,19,2016,2023)
  8 = @(4808.935961081263,601.1169951351578,1165.938996214425,0,3429.901669220172,3924.480503521683,356.7709548656076,489.3711670821389,10.98162164796783,1753.365733305352,25.80253610921124,16.18131984753805,9.386842798111477,64.34086250369539,158,14.36363636363636,21.9055825180374,1,79,2018,2023)
  9 = @(7027.297082973765,878.4121353717206,1175.745185756029,54.36984812050617,3611.18314062814,2115.283199804389,192.2984727094899,337.595798278262,6.098311514417047,1135.983432845721,17.47554070311507,8.949553520309767,6.098311514417047,39.17184251192141,94,8.545454545454545,10.68984224732654,1,30,2016,2023)
  10 = @(8583.02478198114,1430.50413033019,1821.746179999014,1.392243923099459,5037.956787043533,7278.776925880443,661.7069932618584,1479.397583719499,10.98162164796783,5087.076865717208,25.23027567471459,16.04516191724338,10.98162164796783,53.85819572145232,228,20.72727272727273,36.51873740722948,1,128,2017,2023)
  11 = @(4763.591545034287,1587.863848344762,1595.682905637759,629.7445983804416,3429.901669220172,2060.799511978374,686.933170659458,931.6109153745095,31.48722991902208,1753.365733305352,18.05675594648308,3.591709804875617,15.74361495951104,22.1945029532323,98,32.66666666666666,40.8207463593371,2,79,2017,2020)
  12 = @(4763.591545034287,1587.863848344762,1595.682905637759,629.7445983804416,3429.901669220172,2060.799511978374,686.933170659458,931.6109153745095,31.48722991902208,1753.365733305352,18.05675594648308,3.591709804875617,15.74361495951104,22.1945029532323,98,32.66666666666666,40.8207463593371,2,79,2017,2020)
  13 = @(12859.55189450188,918.5394210358487,1239.118252570917,1.392243923099459,5037.956787043533,7801.366065571572,458.9038862100925,1200.333675619392,6.098311514417047,5087.076865717208,22.38754747701466,12.91644130090258,6.098311514417047,56.98420985599093,273,16.05882352941176,29.85688904640622,1,128,2017,2023)
  14 = @(6706.881239126662,2235.627079708887,1419.685520254287,665.9733943173719,3429.901669220172,2415.272480771408,805.0908269238025,862.7300194303008,66.5973394317372,1753.365733305352,33.91111708509359,28.68268333586156,12.94150887031128,66.5973394317372,126,42,39.15354390090378,1,79,2018,2019)
  15 = @(21290.944993777,1330.684062111062,1608.122734911461,54.36984812050617,5037.956787043533,8295.097389970328,460.8387438872404,1164.125099636646,5.258713987810351,5087.076865717208,22.69010555993507,10.50098031613516,5.258713987810351,39.74278801341569,268,14.88888888888889,28.77271765137509,1,128,2016,2023)
  16 = @(28482.6561954059,1294.666190700268,1460.679724284321,0,5037.956787043533,9961.59315927862,415.0663816366092,1052.584154026992,1.087396962410123,5087.076865717208,25.93881637276902,21.11595031539006,1.087396962410123,84.0771316599004,331,13.79166666666667,27.68609628561297,1,128,2016,2023)
  17 = @(6040.90784480929,3020.453922404645,579.0465566296238,2611.006175589118,3429.901669220172,2348.675141339671,1174.337570669835,818.8694805951217,595.3094080343188,1753.365733305352,17.56800591177179,6.542854862312455,12.94150887031128,22.1945029532323,125,62.5,23.33452377915607,46,79,2018,2019)
  18 = @(41403.11895587759,2957.365639705542,5438.578620958389,258.1751543650643,21176.91783511972,8645.091310133475,480.2828505629709,1163.398273955169,6.098311514417047,5087.076865717208,23.69837286756765,14.79328503216368,6.098311514417047,48.38432860277894,282,15.66666666666667,29.08810146879345,1,128,2017,2023)
  19 = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
  20 = @(29361.87809946296,1398.184671402998,1398.856304799021,0,5037.956787043533,9752.995104996171,375.1151963460066,980.436474286025,6.098311514417047,5087.076865717208,18.90658316336147,12.63307169234104,6.098311514417047,55.2973721582137,419,16.11538461538462,27.17841337985266,1,128,2017,2023)
  21 = @(15848.08186278348,990.5051164239674,1131.357598947921,52.14912665617487,3800.513873995446,1553.81569210421,91.40092306495355,144.4937356066764,1.087396962410123,575.827236940745,16.24327655034586,13.03200485417154,1.087396962410123,57.66233027029206,107,6.294117647058823,12.13344914833759,1,50,2017,2023)
  22 = @(2680.405276137013,382.9150394481447,379.1077594348833,0,1088.078154032117,1104.645969257036,122.7384410285595,166.6307284173943,6.098311514417047,538.4389998789497,17.88185964142547,10.90538567613235,6.098311514417047,35.89593332526331,45,5,4.47213595499958,1,15,2018,2022)
}

# Columns G..AA (7..27) of every data row (2..22), in order, per the
# reference dataset that now includes the "meta" statistics.
$colStart = 7

foreach ($rowNum in $rowData.Keys) {
  $values = $rowData[$rowNum]
  for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item([int]$rowNum, $colStart + $i).Value = $values[$i]
  }
}
